$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column headers "Kleur das" -> "Kleur_das" and "Kleur kleding" -> "Kleur_kleding"
# (the decision-tree / sklearn code downstream needs underscore-safe column names)
$ws.Range("F1").Value = "Kleur_das"
$ws.Range("G1").Value = "Kleur_kleding"

# Move the selection, matching the saved workbook state
$ws.Range("G1").Select()
